# Adds a new weekly "Zapallo italiano" price record.
# The new record becomes row 112; the existing rows 112:163 shift down to
# 113:164 (data stays chronologically ordered from newest to oldest).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 112:163 down to 113:164, creating a blank row 112.
$ws.Rows("112:112").Insert()

# Populate the new row 112 with the new weekly record.
$ws.Range("A112").Value = 11
$ws.Range("B112").Value = "Vega Monumental Concepción"
$ws.Range("C112").Value = "Bíobío"
$ws.Range("D112").Value = 44875
$ws.Range("E112").Value = 8
$ws.Range("F112").Value = 100112032
$ws.Range("G112").Value = "Zapallo italiano"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 300
$ws.Range("K112").Value = 7500
$ws.Range("L112").Value = 8000
$ws.Range("M112").Value = 7750
$ws.Range("N112").Value = "$/caja 50 unidades"
$ws.Range("O112").Value = "Región de O'Higgins"
$ws.Range("P112").Value = 155
$ws.Range("Q112").Value = 50
$ws.Range("R112").Value = "Hortaliza"
